# Updated cryptos list on Fri May  5 19:46:54 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) for each coin row,
# and swaps the PEPE / VeChain rows (38 <-> 39) to reflect the new ranking.
#
# For Price cells whose new text would otherwise be auto-parsed by Excel as
# a genuine number (losing formatting like trailing zeros, e.g. "0.07130"),
# the cell is pre-formatted as Text ("@") so the literal string is preserved
# exactly, matching how these values are stored as plain text in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.487.72'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '1.984.82'
$ws.Range("E3").Value = '  +5.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.73'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4682'
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07944'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("E11").Value = '  +4.80%  '
$ws.Range("D12").Value = '1.973.95'
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.247'
$ws.Range("E13").Value = '  +3.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.860'
$ws.Range("E14").Value = '  +3.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07130'
$ws.Range("E15").Value = '  +2.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.64'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009939'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.33'
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '29.589.99'
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.526'
$ws.Range("E22").Value = '  +5.42%  '
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("D24").Value = '2.207.65'
$ws.Range("E24").Value = '  +4.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.105'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.57'
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("E27").Value = '  +1.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.977'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '120.16'
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.958'
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09455'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8936'
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.270'
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.346'
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.176'
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05834'
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.176'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.000003394'
$ws.Range("E38").Value = '  +110.04%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02122'
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.892'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5755'
$ws.Range("E41").Value = '  +1.91%  '
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.810'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.09'
$ws.Range("E44").Value = '  +2.10%  '
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.684'
$ws.Range("E46").Value = '  +6.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.161'
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06945'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.864'
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3090'
$ws.Range("E51").Value = '  +7.72%  '
